# Reposition the Bootstrap Navbar code-snippet screenshots on slides 3-7 so
# that they are no longer covered by the webcam overlay used when recording.
#
# Target Left/Top values below are expressed in points but chosen so that the
# EMU values written back into the OOXML (<a:off x="…" y="…"/>) land exactly
# on the values from the authoritative edit:
#   slide 3 "Picture 5": (6306292,4914906) -> (3814762,5438782)
#   slide 4 "Picture 5": (6750907,4928475) -> (3068140,4826203)
#   slide 4 "Picture 3": (6750907,5669399) -> (3068140,5535176)
#   slide 5 "Picture 6": (5635261,4229973) -> (1642101,5052094)
#   slide 6 "Picture 5": (5635261,4229973) -> (2271275,4138608)
#   slide 6 "Picture 6": (5635261,4949344) -> (2271275,4907399)
#   slide 7 "Picture 4": (5635261,4229973) -> (2288054,4138608)
#   slide 7 "Picture 3": (5635261,4949344) -> (2288054,4903661)
#   slide 7 "Picture 5": (5635261,5640140) -> (2288054,5548775)

$p = $ppt.ActivePresentation

function Move-NamedShape {
    param($slideIndex, $shapeName, $newLeft, $newTop)
    $slide = $p.Slides.Item($slideIndex)
    $shape = $slide.Shapes.Item($shapeName)
    $shape.Left = $newLeft
    $shape.Top = $newTop
}

Move-NamedShape 3 "Picture 5" 300.3749694824219  428.2505798339844
Move-NamedShape 4 "Picture 5" 241.58583068847656 380.0159912109375
Move-NamedShape 4 "Picture 3" 241.58583068847656 435.84063720703125
Move-NamedShape 5 "Picture 6" 129.29930114746094 397.8027038574219
Move-NamedShape 6 "Picture 5" 178.84056091308594 325.8746643066406
Move-NamedShape 6 "Picture 6" 178.84056091308594 386.4093933105469
Move-NamedShape 7 "Picture 4" 180.1617431640625  325.8746643066406
Move-NamedShape 7 "Picture 3" 180.1617431640625  386.11505126953125
Move-NamedShape 7 "Picture 5" 180.1617431640625  436.91143798828125
